$wb = $excel.ActiveWorkbook

# --- Step 1: remove tie rows (OT == ERC, i.e. B == C) from the "OT_wins" sheet ---
# These are languages where OT and ERC produced the same value, so they are not
# really "wins" for either side and are excluded per the commit message
# ("Excluded Lang 14, 19, 59 ... because they're ties").
$wsOT = $wb.Worksheets.Item("OT_wins")

$r = 2
$tieRows = @()
while ($true) {
    $aVal = $wsOT.Cells.Item($r, 1).Value2()
    if ($aVal -eq $null) { break }
    if ($aVal -eq "Average difference") { break }
    $bVal = $wsOT.Cells.Item($r, 2).Value2()
    $cVal = $wsOT.Cells.Item($r, 3).Value2()
    if ($bVal -ne $null -and $cVal -ne $null -and $bVal -eq $cVal) {
        $tieRows += $r
    }
    $r = $r + 1
}

# delete from the bottom up so earlier row numbers stay valid
for ($i = $tieRows.Count - 1; $i -ge 0; $i--) {
    $wsOT.Rows.Item($tieRows[$i]).Delete()
}

# --- Step 2: fix up the "Average difference" formulas on both sheets (drop the $
#     absolute references, and have the OT_wins average start at row 2) ---
$rOT = 2
while ($wsOT.Cells.Item($rOT, 1).Value2() -ne "Average difference") {
    $rOT = $rOT + 1
}
$avgRowOT = $rOT
$lastOT = $avgRowOT - 1
$wsOT.Range("D" + $avgRowOT).Formula = "=AVERAGE(D2:D" + $lastOT + ")"

$wsERC = $wb.Worksheets.Item("ERC_wins")
$rERC = 2
while ($wsERC.Cells.Item($rERC, 1).Value2() -ne "Average difference") {
    $rERC = $rERC + 1
}
$avgRowERC = $rERC
$lastERC = $avgRowERC - 1
$wsERC.Range("D" + $avgRowERC).Formula = "=AVERAGE(D2:D" + $lastERC + ")"

# --- Step 3: rename the sheets to reflect the corrected semantics ---
$wsOT.Name = "ERC_lower_median"
$wsERC.Name = "OT_lower_median"

# --- Step 4: fix up sheet views / selections / active tab ---
$wsERC_lower = $wb.Worksheets.Item("ERC_lower_median")
$wsERC_lower.Select()
$wsERC_lower.Range("D" + $avgRowOT).Select()

$wsOT_lower = $wb.Worksheets.Item("OT_lower_median")
$wsOT_lower.Select()
$wsOT_lower.Range("J10").Select()

$wb.Worksheets.Item("OT_lower_median").Activate()
